# Wish List.xlsx - map renaming / level-table rework
# Sheet1: rows 121-139 hold a small per-level data table. This edit:
#   - removes the three trailing placeholder rows (old A137:A139, levels 16-18)
#   - fills in a header row (row 121, previously blank)
#   - fills in byte-count numbers and per-level notes across the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Drop the three trailing empty levels (old rows 137, 138, 139) ---
$ws.Rows.Item(137).EntireRow.Delete() | Out-Null
$ws.Rows.Item(137).EntireRow.Delete() | Out-Null
$ws.Rows.Item(137).EntireRow.Delete() | Out-Null

# Header row 121 (was entirely blank)
$ws.Cells.Item(121, 2).Value2 = "normal"
$ws.Cells.Item(121, 3).Value2 = "subs"
$ws.Cells.Item(121, 4).Value2 = "reprogram"

# Row 122 (level 1) - unchanged content
# A122 = 1, B122 = "ok" (already present)

# Row 123 (level 2)
$ws.Cells.Item(123, 2).Value2 = 156
$ws.Cells.Item(123, 6).Value2 = "annoying as hell"

# Row 124 (level 3)
$ws.Cells.Item(124, 2).Value2 = 110
$ws.Cells.Item(124, 3).ClearContents() | Out-Null
$ws.Cells.Item(124, 4).Value2 = "62 with reprogram"
$ws.Cells.Item(124, 6).Value2 = "annoying as hell"

# Row 125 (level 4)
$ws.Cells.Item(125, 2).Value2 = 88
$ws.Cells.Item(125, 6).Value2 = "annoying as hell"
$ws.Cells.Item(125, 11).Value2 = "used left switch"

# Row 126 (level 5)
$ws.Cells.Item(126, 2).Value2 = 92
$ws.Cells.Item(126, 6).Value2 = "annoying as hell"
$ws.Cells.Item(126, 11).Value2 = "fixed bad switch"

# Row 127 (level 6)
$ws.Cells.Item(127, 2).Value2 = 78
$ws.Cells.Item(127, 6).Value2 = "very easy, should be an earlier level"

# Row 128 (level 7)
$ws.Cells.Item(128, 2).Value2 = 64
$ws.Cells.Item(128, 6).Value2 = "very easy, should be an earlier level"

# Row 129 (level 8)
$ws.Cells.Item(129, 2).Value2 = 156
$ws.Cells.Item(129, 6).Value2 = "fairly easy - fun water map"

# Row 130 (level 9)
$ws.Cells.Item(130, 2).Value2 = 306
$ws.Cells.Item(130, 3).Value2 = 164
$ws.Cells.Item(130, 6).Value2 = "fun - lots of jumping"
$ws.Cells.Item(130, 8).Value2 = "jump-move forward and jump-move forward 3 are EASILY mass repeated on this map, will test with subs"
$ws.Cells.Item(130, 9).Value2 = "sub1 = jump/move forward, sub2 = jump"

# Row 131 (level 10)
$ws.Cells.Item(131, 2).Value2 = 154
$ws.Cells.Item(131, 6).Value2 = "easy - very straightforward"

# Row 132 (level 11)
$ws.Cells.Item(132, 2).Value2 = 128
$ws.Cells.Item(132, 6).Value2 = "not very complex, fairly easy as well - lots of random extra stuff not related to finishing the map"

# Row 133 (level 12) - unchanged (only A133 = 12)

# Row 134 (level 13)
$ws.Cells.Item(134, 2).Value2 = 182
$ws.Cells.Item(134, 6).Value2 = "could probably lose the reprogram square, additionally, not sure if intended, but the last two switches can be skipped entirely, may want to disable jump on this map, or make some reason for the switches"

# Row 135 (level 14)
$ws.Cells.Item(135, 2).Value2 = 318
$ws.Cells.Item(135, 6).Value2 = "very linear, interesting figuring out what does what"
$ws.Cells.Item(135, 12).Value2 = "needs edge squares removed maybe"

# Row 136 (level 15)
$ws.Cells.Item(136, 6).Value2 = "can't be beaten"

# --- Update the view: scroll down a bit and move the selection to B136 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 118
$ws.Range("B136").Select() | Out-Null
